$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.949.52"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.25"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.68"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4297"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3681"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07263"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.128.40"
$ws.Range("E10").Value = "  +20.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8684"
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.34"
$ws.Range("E12").Value = "  +4.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.414"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.616"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06932"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.13"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008899"
$ws.Range("E18").Value = "  +2.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.28"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.016.51"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.199"
$ws.Range("E22").Value = "  +2.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.02"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.364.84"
$ws.Range("E24").Value = "  +19.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.20"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.885"
$ws.Range("E26").Value = "  +1.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.40"
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.235"
$ws.Range("E28").Value = "  +2.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.922"
$ws.Range("E29").Value = "  +12.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.98"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08972"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.191"
$ws.Range("E32").Value = "  +6.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7501"
$ws.Range("E33").Value = "  +2.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.430"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.807"
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.131"
$ws.Range("E37").Value = "  +4.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05237"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01924"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5116"
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.741"
$ws.Range("E42").Value = "  +8.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.494"
$ws.Range("E43").Value = "  +4.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.361"
$ws.Range("E44").Value = "  +3.50%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "107.11"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.45"
$ws.Range("E46").Value = "  +2.48%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("E47").Value = "  +0.43%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4600"
$ws.Range("E48").Value = "  +2.18%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.649"
$ws.Range("E49").Value = "  +4.08%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06221"
$ws.Range("E50").Value = "  +0.46%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.836"
$ws.Range("E51").Value = "  +4.71%  "
